# Update attendance sheet: set "Invalid" (G) = 1 for row 3, and "Absent" (H) = 1 for rows 3-18
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# G3: Invalid -> 1
$ws.Range("G3").Value = 1

# H3:H18: Absent -> 1
$ws.Range("H3:H18").Value = 1
